# تعديل يدوي في شيت Card21 by admin at 2026-01-14 13:10:33
#
# - B33:K33 were blank placeholder cells; fill them with the literal text "nan"
#   (matching the convention used throughout this sheet for empty numeric fields).
# - Append a new log row (34) for the "قطع سير 1270" event with its own
#   "Serviced by" crew, re-using the same Date/Event/Correction text as row 33.

$wb = $excel.ActiveWorkbook
$ws = $wb.Sheets.Item("Card21")
$ws.Activate()

# Fill the previously-empty measurement columns on row 33 with "nan".
$ws.Range("B33:K33").Value = "nan"

# New row 34: same incident (date/event/correction), different service crew.
$ws.Range("L34").Value = "14/1/2026"
$ws.Range("M34").Value = "قطع سير 1270"
$ws.Range("N34").Value = "تم تغير سير 1270(مشلان)"
$ws.Range("O34").Value = "محمود ايهاب،سعيد،م.محمد عبدالله "
